$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.219.98"
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").Value = "3.407.80"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.71"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "619.08"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.48%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.43"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +3.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.390"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.975"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +3.45%  "

$ws.Range("D11").Value = "3.409.07"
$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.08"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +6.57%  "

$ws.Range("E13").Value = "  +1.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.27"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +4.18%  "

$ws.Range("D15").Value = "93.057.72"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").Value = "4.056.59"
$ws.Range("E16").Value = "  +1.37%  "

$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.22"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +2.68%  "

$ws.Range("D19").Value = "3.398.34"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.94"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +6.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.65"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +6.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.499"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +10.47%  "

$ws.Range("E23").Value = "  +7.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "496.32"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.72"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +6.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000182"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -2.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "90.52"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.99"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +4.18%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.28"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.138"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +4.88%  "

$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("E34").Value = "  +1.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.547"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +3.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.88"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "558.78"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +6.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.44"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.68%  "

$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.06%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.40"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("E41").Value = "  +1.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.895"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +1.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.68"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -1.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.65"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.63%  "

$ws.Range("E45").Value = "  +1.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0412"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +4.51%  "

$ws.Range("E47").Value = "  -1.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.92"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.10"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -2.30%  "

$ws.Range("E50").Value = "  +1.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.06"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.76%  "
